$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged, values updated)
$ws.Range("B3").Value = 0.03730215356940209
$ws.Range("C3").Value = 0.03715467150880781
$ws.Range("D3").Value = 0.06174247312753128

# Row 4 - renamed GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03402465816345154
$ws.Range("C4").Value = 0.03514819938377833
$ws.Range("D4").Value = 0.09199082799209106

# Row 5 - renamed AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.04501843076036613
$ws.Range("C5").Value = 0.03271669992482031
$ws.Range("D5").Value = 0.03525352047568234
